{"js": "// The target document is a 5-row x 3-column table of \"lattice multiplication\"\n// practice cells. Each cell holds a single paragraph / single run with text\n// like:\n//   \"94 x 21\" <br/> \"  2    1\" <br/> \"  ----\" <br/> \"9|    |\" <br/> \"4|    |\"\n// (the line breaks are manual <w:br/> line breaks, which Word's object\n// model exposes as vertical-tab (\\u000b) characters inside Range/Paragraph\n// .text). The edit simply swaps every cell's five lines for a new set of\n// multiplication problems, while leaving the 5x3 table shape, the run\n// formatting (sz=32) and everything else untouched.\n\n// New content for every cell, in row-major order (row 0 cell 0, row 0 cell\n// 1, row 0 cell 2, row 1 cell 0, ...). Each entry is the 5 lines joined with\n// the vertical-tab line-break character Word uses for manual breaks.\nconst NEW_CELL_TEXT = [\n  \"94 x 21\\u000b  2    1\\u000b  ----\\u000b9|    |\\u000b4|    |\",\n  \"57 x 70\\u000b  7    0\\u000b  ----\\u000b5|    |\\u000b7|    |\",\n  \"23 x 69\\u000b  6    9\\u000b  ----\\u000b2|    |\\u000b3|    |\",\n  \"98 x 84\\u000b  8    4\\u000b  ----\\u000b9|    |\\u000b8|    |\",\n  \"79 x 63\\u000b  6    3\\u000b  ----\\u000b7|    |\\u000b9|    |\",\n  \"60 x 82\\u000b  8    2\\u000b  ----\\u000b6|    |\\u000b0|    |\",\n  \"32 x 95\\u000b  9    5\\u000b  ----\\u000b3|    |\\u000b2|    |\",\n  \"60 x 60\\u000b  6    0\\u000b  ----\\u000b6|    |\\u000b0|    |\",\n  \"63 x 72\\u000b  7    2\\u000b  ----\\u000b6|    |\\u000b3|    |\",\n  \"64 x 24\\u000b  2    4\\u000b  ----\\u000b6|    |\\u000b4|    |\",\n  \"76 x 50\\u000b  5    0\\u000b  ----\\u000b7|    |\\u000b6|    |\",\n  \"11 x 52\\u000b  5    2\\u000b  ----\\u000b1|    |\\u000b1|    |\",\n  \"28 x 52\\u000b  5    2\\u000b  ----\\u000b2|    |\\u000b8|    |\",\n  \"27 x 63\\u000b  6    3\\u000b  ----\\u000b2|    |\\u000b7|    |\",\n  \"81 x 31\\u000b  3    1\\u000b  ----\\u000b8|    |\\u000b1|    |\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\n// Load the cells of every row up front.\nfor (const row of table.rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Load the single paragraph inside every cell (each cell body is one\n// paragraph / one run in this document).\nfor (const row of table.rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.paragraphs.load(\"items\");\n  }\n}\nawait context.sync();\n\nlet idx = 0;\nfor (const row of table.rows.items) {\n  for (const cell of row.cells.items) {\n    const para = cell.body.paragraphs.items[0];\n    // getRange(\"Content\") is the paragraph's text range *without* the\n    // trailing paragraph mark, so replacing it keeps the existing run's\n    // formatting (sz=32) instead of Word re-deriving/duplicating it onto a\n    // new paragraph-mark run.\n    const contentRange = para.getRange(\"Content\");\n    contentRange.insertText(NEW_CELL_TEXT[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The target document is a 5-row x 3-column table of \"lattice multiplication\"\n# practice cells. Each cell holds a single paragraph / single run with text\n# like:\n#   \"94 x 21\" <br/> \"  2    1\" <br/> \"  ----\" <br/> \"9|    |\" <br/> \"4|    |\"\n# (the line breaks are manual line breaks -- Word represents a manual\n# line break as a vertical-tab character, chr(11)/0x0B, inside Range.Text).\n# The edit simply swaps every cell's five lines for a new set of\n# multiplication problems, while leaving the 5x3 table shape, the run\n# formatting (sz=32) and everything else untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$lineBreak = [char]11\n\n# New content for every cell, addressed by (row, column) -- 1-based, matching\n# Word's Cell(row, column) indexing. Each entry is the 5 lines for that cell.\n$newCells = @{\n    \"1,1\" = @(\"94 x 21\", \"  2    1\", \"  ----\", \"9|    |\", \"4|    |\")\n    \"1,2\" = @(\"57 x 70\", \"  7    0\", \"  ----\", \"5|    |\", \"7|    |\")\n    \"1,3\" = @(\"23 x 69\", \"  6    9\", \"  ----\", \"2|    |\", \"3|    |\")\n    \"2,1\" = @(\"98 x 84\", \"  8    4\", \"  ----\", \"9|    |\", \"8|    |\")\n    \"2,2\" = @(\"79 x 63\", \"  6    3\", \"  ----\", \"7|    |\", \"9|    |\")\n    \"2,3\" = @(\"60 x 82\", \"  8    2\", \"  ----\", \"6|    |\", \"0|    |\")\n    \"3,1\" = @(\"32 x 95\", \"  9    5\", \"  ----\", \"3|    |\", \"2|    |\")\n    \"3,2\" = @(\"60 x 60\", \"  6    0\", \"  ----\", \"6|    |\", \"0|    |\")\n    \"3,3\" = @(\"63 x 72\", \"  7    2\", \"  ----\", \"6|    |\", \"3|    |\")\n    \"4,1\" = @(\"64 x 24\", \"  2    4\", \"  ----\", \"6|    |\", \"4|    |\")\n    \"4,2\" = @(\"76 x 50\", \"  5    0\", \"  ----\", \"7|    |\", \"6|    |\")\n    \"4,3\" = @(\"11 x 52\", \"  5    2\", \"  ----\", \"1|    |\", \"1|    |\")\n    \"5,1\" = @(\"28 x 52\", \"  5    2\", \"  ----\", \"2|    |\", \"8|    |\")\n    \"5,2\" = @(\"27 x 63\", \"  6    3\", \"  ----\", \"2|    |\", \"7|    |\")\n    \"5,3\" = @(\"81 x 31\", \"  3    1\", \"  ----\", \"8|    |\", \"1|    |\")\n}\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $key = \"$r,$c\"\n        $lines = $newCells[$key]\n        $cell = $t.Cell($r, $c)\n        # Setting Range.Text directly (rather than e.g. deleting then\n        # inserting a new paragraph) keeps the existing run's formatting\n        # (sz=32) and avoids touching the cell's end-of-cell mark.\n        $cell.Range.Text = [string]::Join($lineBreak, $lines)\n    }\n}\n"}
